# Fixed naive component forecaster bug - Presentation state 11.02.
# For rows 2..16, shift the existing values in B:J one column to the right
# (into C:K), dropping whatever previously sat in column K, and insert a
# newly computed naive-forecast value into column B.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newB = @{
  2  = -0.2177157015159319
  3  = -0.1395947820665385
  4  = -0.3119065001142551
  5  = 0.7021231295320197
  6  = 1.514070997382048
  7  = 0.2163102553365951
  8  = 0.3684555432821496
  9  = 0.661541622456546
  10 = -0.07992401592518952
  11 = 0.1551026493581833
  12 = -0.08373363042288225
  13 = 0.1925427069667326
  14 = -0.4379379024501944
  15 = 0.2324016585002178
  16 = -0.09587373626955231
}

for ($r = 2; $r -le 16; $r++) {
    # Shift existing values right, starting from the rightmost populated
    # column so values are not overwritten before being moved.
    for ($c = 10; $c -ge 2; $c--) {
        $srcCell = $ws.Cells.Item($r, $c)
        $srcValue = $srcCell.Value()
        if ($srcValue -ne $null) {
            $ws.Cells.Item($r, $c + 1).Value = $srcValue
        }
    }
    # Insert the new leading value into column B.
    $ws.Cells.Item($r, 2).Value = $newB[$r]
}
